$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Case Matrix")

# Rows that get the text "No discount given" in column I (Expected Result)
$textRows = @(19,20,21,23,24,25,26,27,28,29,30,31,32,33,35,37,39,41,42)
foreach ($r in $textRows) {
    $ws.Cells.Item($r, 9).Value = "No discount given"
}

# Rows that get a numeric percentage discount value in column I
$percentRows = @{
    22 = 0.1
    34 = 0.15
    36 = 0.25
    40 = 0.15
}
foreach ($r in $percentRows.Keys) {
    $cell = $ws.Cells.Item($r, 9)
    $cell.Value = $percentRows[$r]
    $cell.NumberFormat = "0%"
}

# Update selection to reflect where the author ended up editing
$ws.Activate()
[void]$ws.Range("I43").Select()
